# Updated Timesheet for Week 6: logged a new entry on 2/14 for
# "Created more pages for the site" (12:30 PM - 2:45 PM, 2.25 hrs).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 6")
$ws.Activate()

# Row 4 sits in the previously-blank gap between the last logged day (row 3)
# and the weekly-total rows (20/21), so this is a plain fill-in, not a row
# shift. Seed its formatting from row 2 (same column layout: date / start /
# end / description / hours) before writing the new values.
$ws.Range("A2:E2").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 43509
$ws.Range("B4").Value = 0.52083333333333337
$ws.Range("C4").Value = 0.61458333333333337
$ws.Range("D4").Value = "Created more pages for the site"
$ws.Range("E4").Value = 2.25

$ws.Rows.Item(4).RowHeight = 18

$ws.Range("D4").Select()
